$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 updates (picks up the values previously on row 16)
$ws.Range("A14").Value = 111798755
$ws.Range("Q14").Value = 753030.7189070459
$ws.Range("R14").Value = 7090920.781295684
$ws.Range("S14").Value = 25
# AR14 becomes a present-but-empty text cell (mirrors row 16's old AR cell).
# A plain "" assignment deletes the cell outright, so write a lone quote
# (Excel's text-prefix marker) then strip the formatting it introduces,
# leaving a real empty string value in place.
$ws.Range("AR14").Value = "'"
$ws.Range("AR14").ClearFormats()

# Row 15 updates
$ws.Range("A15").Value = 111798795
$ws.Range("AF15").Value = "'"
$ws.Range("AF15").ClearFormats()
$ws.Range("AI15").ClearContents()
$ws.Range("AR15").ClearContents()

# Row 16 updates (picks up the values previously on row 14)
$ws.Range("A16").Value = 111798760
$ws.Range("Q16").Value = 753108.8301749222
$ws.Range("R16").Value = 7091007.708399305
$ws.Range("S16").Value = 100
$ws.Range("AR16").ClearContents()
